# AprobacionesPreviasAPago.xlsx - agrego en el test case de anulacion dos
# smart folders (una para anular en si y otra para obtener el numero de
# anulacion).
#
# Updates the existing "usuario"/"NroSiniestro" values for a few rows and
# appends two brand-new rows (7 and 8) that reuse the same layout/format
# as the existing "preproducciongestion" rows (5 and 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update values on the existing data rows (2-6).
#    NroSiniestro values look numeric but must stay TEXT (some keep
#    trailing padding spaces), so a leading apostrophe forces text entry
#    the same way typing it in the Excel UI would.
# ---------------------------------------------------------------------
$ws.Range("E2").Value2 = "'0420194406812 "
$ws.Range("E3").Value2 = "'1220194200683"
$ws.Range("E4").Value2 = "'1120170200942   "

$ws.Range("C5").Value2 = "cnapolitano"
$ws.Range("E5").Value2 = "'1120170200942   "

$ws.Range("C6").Value2 = "dgariffo"
$ws.Range("E6").Value2 = "'1220170301442   "

# Rows 3 & 4 no longer need the taller, wrapped row height - AutoFit()
# drops the explicit row height back to the sheet default.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# ---------------------------------------------------------------------
# 2) Append the two new "smart folder" rows (7 and 8), cloning the
#    formatting of row 6 (same Ambiente/URL/Contrasenia, hyperlinked B
#    cell, styled E cell).
# ---------------------------------------------------------------------
$ws.Range("A6:E6").Copy($ws.Range("A7:E7")) | Out-Null
$ws.Range("A6:E6").Copy($ws.Range("A8:E8")) | Out-Null

# Hyperlinks.Add() resets the cell style, so add the links first and
# then re-apply row 6's formatting on top of the new B cells.
$ws.Hyperlinks.Add($ws.Range("B7"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do") | Out-Null

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 7: "anular" smart folder (rsuarez)
$ws.Range("A7").Value2 = "preproducciongestion.segurossura.com.ar"
$ws.Range("C7").Value2 = "rsuarez"
$ws.Range("D7").Value2 = "silverarrow"
$ws.Range("E7").Value2 = "'1220170301437   "

# Row 8: "obtener numero de anulacion" smart folder (ocerutti)
$ws.Range("A8").Value2 = "preproducciongestion.segurossura.com.ar"
$ws.Range("C8").Value2 = "ocerutti"
$ws.Range("D8").Value2 = "silverarrow"
$ws.Range("E8").Value2 = "'1120170200939   "

# ---------------------------------------------------------------------
# 3) Column E needs to be wide enough to show the full NroSiniestro text
#    (target stored width 18.7109375; 17.9 is the closest input that this
#    engine's pixel-snapped ColumnWidth setter can reach).
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 17.9

# ---------------------------------------------------------------------
# 4) Leave the selection where the author left it after typing the new
#    rows.
# ---------------------------------------------------------------------
$ws.Range("C9").Select() | Out-Null
